$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 4 (Operating System / Mr Deepak),
# pushing it down to row 5, and creating space for the new EVS / Mr Tony class.
$ws.Rows("4:4").Insert()

# Fill in the new class row (times copy the 0.75 / 0.83333333333333337 fractions
# Excel uses for 18:00 and 20:00).
$ws.Range("A4").Value = 0.75
$ws.Range("B4").Value = 0.83333333333333337
$ws.Range("C4").Value = "EVS"
$ws.Range("D4").Value = "Mr Tony"

# New "email" column header and the teacher's e-mail address (with hyperlink).
$ws.Range("E1").Value = "email"
$ws.Range("E4").Value = "girishhrudhay@gmail.com"
$ws.Hyperlinks.Add($ws.Range("E4"), "mailto:girishhrudhay@gmail.com", "", "", "girishhrudhay@gmail.com")

$ws.Range("J10").Select()
